# Generate Report for Handback
#
# The localization-status report is regenerated: the "b3b73927-..." row has
# been handed back (target file is now in sync with en-US), so:
#   - the Status column flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is shown
#     (Overview!E2/F2/E3/F3 and the zh-cn/de-de "Status" column)
#   - the zh-cn and de-de per-language sheets gain a "Latest Target File"
#     hyperlink (back to the source .md) and a "Latest Handback File" /
#     "Latest Handback DateTime" for that row

$wb = $excel.ActiveWorkbook

$mdFile      = "b3b73927-7c3c-4aeb-a33b-5ad6b1e94038.md"
$mdUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e094793d28943dcdb2a38416208554848878e71/e2e/$mdFile"
$statusText  = "Handed back: in sync with en-US"
$zhHandback  = "b3b73927-7c3c-4aeb-a33b-5ad6b1e94038.2d9e3e063bc9b5c1ad485af82df07b1092b00948.zh-cn.xlf"
$deHandback  = "b3b73927-7c3c-4aeb-a33b-5ad6b1e94038.2d9e3e063bc9b5c1ad485af82df07b1092b00948.de-de.xlf"
$zhDateTime  = "2016-09-02 23:09:46"
$deDateTime  = "2016-09-02 23:09:53"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status cells for the first file (rows 2,3
# both reference the same file in this workbook) go from "Ready for
# handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Latest Target File -> hyperlink back to the source markdown file
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, "", "", $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl, "", "", $mdFile)

# Latest Handback File / DateTime
$wsZh.Range("J2").Value = $zhHandback
$wsZh.Range("J3").Value = $zhHandback
$wsZh.Range("K2").Value = $zhDateTime
$wsZh.Range("K3").Value = $zhDateTime

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Latest Target File -> hyperlink back to the source markdown file
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, "", "", $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl, "", "", $mdFile)

# Latest Handback File / DateTime
$wsDe.Range("J2").Value = $deHandback
$wsDe.Range("J3").Value = $deHandback
$wsDe.Range("K2").Value = $deDateTime
$wsDe.Range("K3").Value = $deDateTime

# ---------------------------------------------------------------------
# Widen the columns that now hold the longer strings, matching the
# regenerated report's auto-fit widths.
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
